# "Conclusao da CadastroPage referente ao cadastro usando massa do Excel"
# Update the registration ("cadastro") test-data row used by the Excel-driven
# data provider: swap user1 -> user4 for the plan/city/address fields.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values: username -> user4Plan, city -> cidadePlan, address -> "Rua Test Plan, 1 "
$ws.Range("A2").Value = "user4Plan"
$ws.Range("H2").Value = "cidadePlan"
$ws.Range("I2").Value = "Rua Test Plan, 1 "

# Widen the address column (I) so the longer text fits, same as Excel does
# automatically when a cell's content no longer fits the column.
$ws.Columns.Item(9).ColumnWidth = 16.3

# Active selection ends up on D13 after the edit (e.g. from scrolling down
# while reviewing the sheet).
$ws.Range("D13").Select()
